$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update text labels ---
$ws.Range("F2").Value = "4 minutes and 46.5 seconds"
$ws.Range("C7").Value = "add the 30s"
$ws.Range("A8").Value = "After the 3 first codes and the 30s wait"

# --- New label cells for the "time to enter a code" calculation ---
# (order matters for shared-string table layout: I1, then H1, then L1)
$ws.Range("I1").Value = "time for each 3 digit code to be entered"
$ws.Range("H1").Value = "delay time between button presses (s)"
$ws.Range("L1").Value = "pressing enter time in seconds"

# --- Length of timeout is now computed instead of a hard-coded value ---
$ws.Range("E2").Formula = "=4.775*60"

# --- New supporting values/formulas ---
$ws.Range("H2").Value = 0.114
$ws.Range("L2").Formula = "=H2+0.07"
$ws.Range("I2").Formula = "=B1*(2*H2)+L2"

# --- Worst case time now also accounts for time spent entering codes, and timeout wait changed from 20s to 30s ---
$ws.Range("B8").Formula = "=(((B5/E1)*E2)+((B2)*I2))"
$ws.Range("C8").Formula = "=B8+30"
